$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 444.45715
$ws.Range("J19").Value = 486.85
$ws.Range("L19").Value = 486.85
$ws.Range("N19").Value = -836.85
$ws.Range("H28").Value = 505920.72
$ws.Range("I28").Value = 694894.1
$ws.Range("J28").Value = 1991.6666
$ws.Range("K28").Value = 694894.1
$ws.Range("L28").Value = 1991.6666
$ws.Range("M28").Value = -694409.1
$ws.Range("N28").Value = -2961.6666
$ws.Range("H107").Value = 654016.5
$ws.Range("I107").Value = 654016.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 654016.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -652096.5
$ws.Range("N107").ClearContents()
$ws.Range("H111").Value = 1408.7142
$ws.Range("I111").Value = 880
$ws.Range("J111").Value = 1702.4445
$ws.Range("K111").Value = 2640
$ws.Range("L111").Value = 5107.333500000001
$ws.Range("M111").Value = 427
$ws.Range("N111").Value = -11241.3335
$ws.Range("H113").Value = 5696.6
$ws.Range("I113").Value = 3790
$ws.Range("J113").Value = 6967.6665
$ws.Range("K113").Value = 3790
$ws.Range("L113").Value = 6967.6665
$ws.Range("M113").Value = -536
$ws.Range("N113").Value = -13475.6665
$ws.Range("H132").Value = 22902.447
$ws.Range("I132").Value = 25266.738
$ws.Range("J132").Value = 3042.4
$ws.Range("K132").Value = 75800.21400000001
$ws.Range("L132").Value = 9127.200000000001
$ws.Range("M132").Value = -73270.21400000001
$ws.Range("N132").Value = -14187.2
$ws.Range("H138").Value = 9075483
$ws.Range("I138").Value = 2317514.2
$ws.Range("J138").Value = 14288772
$ws.Range("K138").Value = 6952542.600000001
$ws.Range("L138").Value = 42866316
$ws.Range("M138").Value = -6947402.600000001
$ws.Range("N138").Value = -42876596

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 145747.58
$ws.Range("I2").Value = 145747.58
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 145747.58
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -145634.58
$ws.Range("N2").ClearContents()
$ws.Range("H61").Value = 2081.9119
$ws.Range("I61").Value = 1334.4445
$ws.Range("J61").Value = 4965
$ws.Range("K61").Value = 1334.4445
$ws.Range("L61").Value = 4965
$ws.Range("M61").Value = -1122.4445
$ws.Range("N61").Value = -5389
$ws.Range("H104").Value = 36333.332
$ws.Range("J104").Value = 36333.332
$ws.Range("L104").Value = 36333.332
$ws.Range("N104").Value = -43321.332
$ws.Range("H116").Value = 145747.58
$ws.Range("I116").Value = 145747.58
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 145747.58
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -143453.58
$ws.Range("N116").ClearContents()
$ws.Range("H136").Value = 2081.9119
$ws.Range("I136").Value = 1334.4445
$ws.Range("J136").Value = 4965
$ws.Range("K136").Value = 4003.3335
$ws.Range("L136").Value = 14895
$ws.Range("M136").Value = -1453.3335
$ws.Range("N136").Value = -19995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 145747.58
$ws.Range("I3").Value = 145747.58
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 145747.58
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -145633.58
$ws.Range("N3").ClearContents()
$ws.Range("H99").Value = 1617
$ws.Range("I99").Value = 1330.7273
$ws.Range("K99").Value = 1330.7273
$ws.Range("M99").Value = 167.2727
$ws.Range("H103").Value = 27714.285
$ws.Range("J103").Value = 27714.285
$ws.Range("L103").Value = 27714.285
$ws.Range("N103").Value = -30058.285
$ws.Range("H107").Value = 2499.6667
$ws.Range("I107").Value = 2499.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2499.6667
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -579.6667000000002
$ws.Range("N107").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1721.0741
$ws.Range("I58").Value = 1166.3914
$ws.Range("J58").Value = 4910.5
$ws.Range("K58").Value = 1166.3914
$ws.Range("L58").Value = 4910.5
$ws.Range("M58").Value = -963.3914
$ws.Range("N58").Value = -5316.5
$ws.Range("H107").Value = 490.77777
$ws.Range("I107").Value = 490.77777
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 490.77777
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1429.22223
$ws.Range("N107").ClearContents()
$ws.Range("H136").Value = 1721.0741
$ws.Range("I136").Value = 1166.3914
$ws.Range("J136").Value = 4910.5
$ws.Range("K136").Value = 3499.1742
$ws.Range("L136").Value = 14731.5
$ws.Range("M136").Value = -949.1741999999999
$ws.Range("N136").Value = -19831.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1676.9
$ws.Range("I99").Value = 1658.625
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 4975.875
$ws.Range("L99").Value = 5250
$ws.Range("M99").Value = -2729.875
$ws.Range("N99").Value = -9742
$ws.Range("H113").Value = 777
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 777
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2331
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6671
$ws.Range("H131").Value = 4903266.5
$ws.Range("I131").Value = 327.75
$ws.Range("J131").Value = 5953896
$ws.Range("K131").Value = 983.25
$ws.Range("L131").Value = 17861688
$ws.Range("M131").Value = 4056.75
$ws.Range("N131").Value = -17871768

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1190.4615
$ws.Range("I107").Value = 906.9091
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 906.9091
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = 1013.0909
$ws.Range("N107").Value = -6590
$ws.Range("H136").Value = 18904.625
$ws.Range("J136").Value = 18683.467
$ws.Range("L136").Value = 56050.401
$ws.Range("N136").Value = -61150.401

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 937.0625
$ws.Range("I61").Value = 588.1111
$ws.Range("J61").Value = 1385.7142
$ws.Range("K61").Value = 588.1111
$ws.Range("L61").Value = 1385.7142
$ws.Range("M61").Value = -386.1111
$ws.Range("N61").Value = -1789.7142
$ws.Range("H104").Value = 26360
$ws.Range("J104").Value = 26360
$ws.Range("L104").Value = 26360
$ws.Range("N104").Value = -33348
$ws.Range("H106").Value = 26000
$ws.Range("J106").Value = 26000
$ws.Range("L106").Value = 26000
$ws.Range("N106").Value = -28524
$ws.Range("H113").Value = 937.0625
$ws.Range("I113").Value = 588.1111
$ws.Range("J113").Value = 1385.7142
$ws.Range("K113").Value = 588.1111
$ws.Range("L113").Value = 1385.7142
$ws.Range("M113").Value = 1581.8889
$ws.Range("N113").Value = -5725.7142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H107").Value = 1326
$ws.Range("I107").Value = 1326
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3978
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2058
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 2916.8696
$ws.Range("I132").Value = 2848.4412
$ws.Range("J132").Value = 3110.75
$ws.Range("K132").Value = 8545.3236
$ws.Range("L132").Value = 9332.25
$ws.Range("M132").Value = -6015.3236
$ws.Range("N132").Value = -14392.25
$ws.Range("H136").Value = 2946.4807
$ws.Range("I136").Value = 1137.4865
$ws.Range("J136").Value = 7408.6665
$ws.Range("K136").Value = 3412.4595
$ws.Range("L136").Value = 22225.9995
$ws.Range("M136").Value = -862.4594999999999
$ws.Range("N136").Value = -27325.9995
